# Applies the "Updated cryptos list" data refresh to sheet1 (cryptos.xlsx).
# Row-for-row cell updates (Price / Volume(1h) columns, and a few reordered
# Coin/Link/Price/Volume rows where the row assignment of coins changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.238.43"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.656.28"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07726"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.599"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "1.622.23"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "1.883.94"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "0.0₅8256"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "26.236.58"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.698"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1203"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.275"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.514"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05638"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.279"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9548"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.803"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5760"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8436"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.010.08"
$ws.Range("E44").Value = "  -6.00%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.794.94"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05351"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈104"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4349"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.016"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
